# "Generate Report for Handback" - refresh the localization-status report
# after a handback: status flips from "Ready for handoff" to
# "Handed back: in sync with en-US", the handback timestamps advance, the
# stale "handback not latest" error is cleared, and Excel auto-fits the
# Status / Error Detail columns to their new (longer / shorter) content.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---- Overview sheet: per-locale status columns (E = zh-cn, F = de-de) ----
# These mirror the same "Status" text shown on the locale sheets, so they
# roll over to the new status too, and widen to fit it.
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E2").Value = $newStatus
$ov.Range("F2").Value = $newStatus
$ov.Range("E3").Value = $newStatus
$ov.Range("F3").Value = $newStatus
$ov.Columns.Item(5).ColumnWidth = 29.15
$ov.Columns.Item(6).ColumnWidth = 29.15

# ---- zh-cn sheet ----
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Columns.Item(3).ColumnWidth = 29.15   # Status column, now holds a longer string
$zh.Columns.Item(16).ColumnWidth = 12.8   # Error Detail column, now empty

$zh.Range("C2").Value = $newStatus
$zh.Range("C3").Value = $newStatus

$zh.Range("K2").Value = "2016-11-02 05:17:49"
$zh.Range("K3").Value = "2016-11-02 05:17:49"

$zh.Range("P2").Value = ""

# ---- de-de sheet ----
$de = $wb.Worksheets.Item("de-de")
$de.Columns.Item(3).ColumnWidth = 29.15   # Status column, now holds a longer string
$de.Columns.Item(16).ColumnWidth = 12.8   # Error Detail column, now empty

$de.Range("C2").Value = $newStatus
$de.Range("C3").Value = $newStatus

$de.Range("K2").Value = "2016-11-02 05:18:08"
$de.Range("K3").Value = "2016-11-02 05:18:08"

$de.Range("P2").Value = ""
